# Scheduled market-data refresh: updates Universalis price snapshots
# (currentAveragePrice / NQ / HQ) and the derived Leve-turn-in profit
# columns (H:N) on each job sheet of the Bahamut_Profits workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 34500
$ws.Range("I13").Value = 34000
$ws.Range("K13").Value = 34000
$ws.Range("M13").Value = -33831
$ws.Range("H76").Value = 50002780
$ws.Range("I76").Value = 53574200
$ws.Range("J76").Value = 2900
$ws.Range("K76").Value = 53574200
$ws.Range("L76").Value = 2900
$ws.Range("M76").Value = -53573885
$ws.Range("N76").Value = -3530
$ws.Range("H79").Value = 50002780
$ws.Range("I79").Value = 53574200
$ws.Range("J79").Value = 2900
$ws.Range("K79").Value = 53574200
$ws.Range("L79").Value = 2900
$ws.Range("M79").Value = -53573108
$ws.Range("N79").Value = -5084
$ws.Range("H100").Value = 918.2
$ws.Range("I100").Value = 1023.75
$ws.Range("J100").Value = 496
$ws.Range("K100").Value = 1023.75
$ws.Range("L100").Value = 496
$ws.Range("M100").Value = -482.75
$ws.Range("N100").Value = -1578
$ws.Range("H129").Value = 1950265
$ws.Range("J129").Value = 3087659.8
$ws.Range("L129").Value = 9262979.399999999
$ws.Range("N129").Value = -9272979.399999999
$ws.Range("H138").Value = 4267.8247
$ws.Range("I138").Value = 3149.077
$ws.Range("J138").Value = 4440.9644
$ws.Range("K138").Value = 9447.231
$ws.Range("L138").Value = 13322.8932
$ws.Range("M138").Value = -4307.231
$ws.Range("N138").Value = -23602.8932

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13128.12
$ws.Range("I32").Value = 11230.901
$ws.Range("J32").Value = 32311.111
$ws.Range("K32").Value = 11230.901
$ws.Range("L32").Value = 32311.111
$ws.Range("M32").Value = -10943.901
$ws.Range("N32").Value = -32885.111
$ws.Range("H45").Value = 1565.0454
$ws.Range("I45").Value = 1819.7273
$ws.Range("J45").Value = 1310.3636
$ws.Range("K45").Value = 1819.7273
$ws.Range("L45").Value = 1310.3636
$ws.Range("M45").Value = -1442.7273
$ws.Range("N45").Value = -2064.3636
$ws.Range("H102").Value = 2586.3333
$ws.Range("I102").Value = 3377.25
$ws.Range("J102").Value = 1004.5
$ws.Range("K102").Value = 3377.25
$ws.Range("L102").Value = 1004.5
$ws.Range("M102").Value = -1755.25
$ws.Range("N102").Value = -4248.5
$ws.Range("H132").Value = 2315.3225
$ws.Range("I132").Value = 1571.3334
$ws.Range("J132").Value = 3345.4614
$ws.Range("K132").Value = 4714.0002
$ws.Range("L132").Value = 10036.3842
$ws.Range("M132").Value = -2184.0002
$ws.Range("N132").Value = -15096.3842

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2722.6223
$ws.Range("I105").Value = 2841.861
$ws.Range("J105").Value = 2245.6667
$ws.Range("K105").Value = 2841.861
$ws.Range("L105").Value = 2245.6667
$ws.Range("M105").Value = -1094.861
$ws.Range("N105").Value = -5739.6667
$ws.Range("H134").Value = 12704.034
$ws.Range("I134").Value = 1138.3422
$ws.Range("J134").Value = 85953.414
$ws.Range("K134").Value = 3415.0266
$ws.Range("L134").Value = 257860.242
$ws.Range("M134").Value = -880.0266000000001
$ws.Range("N134").Value = -262930.242

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3106.025
$ws.Range("I31").Value = 2578.2964
$ws.Range("J31").Value = 4202.077
$ws.Range("K31").Value = 2578.2964
$ws.Range("L31").Value = 4202.077
$ws.Range("M31").Value = -2283.2964
$ws.Range("N31").Value = -4792.077
$ws.Range("H34").Value = 3106.025
$ws.Range("I34").Value = 2578.2964
$ws.Range("J34").Value = 4202.077
$ws.Range("K34").Value = 2578.2964
$ws.Range("L34").Value = 4202.077
$ws.Range("M34").Value = -2376.2964
$ws.Range("N34").Value = -4606.077
$ws.Range("H58").Value = 3156.9111
$ws.Range("I58").Value = 802.21875
$ws.Range("J58").Value = 8953.076999999999
$ws.Range("K58").Value = 802.21875
$ws.Range("L58").Value = 8953.076999999999
$ws.Range("M58").Value = -599.21875
$ws.Range("N58").Value = -9359.076999999999
$ws.Range("H62").Value = 6709.4443
$ws.Range("I62").Value = 6897.5
$ws.Range("K62").Value = 6897.5
$ws.Range("M62").Value = -6273.5
$ws.Range("H65").Value = 6709.4443
$ws.Range("I65").Value = 6897.5
$ws.Range("K65").Value = 34487.5
$ws.Range("M65").Value = -31367.5
$ws.Range("H132").Value = 1521.4857
$ws.Range("I132").Value = 989.1667
$ws.Range("J132").Value = 2682.9092
$ws.Range("K132").Value = 2967.5001
$ws.Range("L132").Value = 8048.7276
$ws.Range("M132").Value = -437.5001000000002
$ws.Range("N132").Value = -13108.7276
$ws.Range("H136").Value = 3156.9111
$ws.Range("I136").Value = 802.21875
$ws.Range("J136").Value = 8953.076999999999
$ws.Range("K136").Value = 2406.65625
$ws.Range("L136").Value = 26859.231
$ws.Range("M136").Value = 143.34375
$ws.Range("N136").Value = -31959.231

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1537.1034
$ws.Range("I122").Value = 599
$ws.Range("J122").Value = 1645.3462
$ws.Range("K122").Value = 5391
$ws.Range("L122").Value = 14808.1158
$ws.Range("M122").Value = -2941
$ws.Range("N122").Value = -19708.1158

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11669.77
$ws.Range("I70").Value = 14702.135
$ws.Range("K70").Value = 14702.135
$ws.Range("M70").Value = -14432.135
$ws.Range("H73").Value = 11669.77
$ws.Range("I73").Value = 14702.135
$ws.Range("K73").Value = 14702.135
$ws.Range("M73").Value = -13766.135
$ws.Range("H97").Value = 2426.1875
$ws.Range("I97").Value = 2522
$ws.Range("J97").Value = 2266.5
$ws.Range("K97").Value = 2522
$ws.Range("L97").Value = 2266.5
$ws.Range("M97").Value = -2026
$ws.Range("N97").Value = -3258.5
$ws.Range("H109").Value = 30213.75
$ws.Range("J109").Value = 30213.75
$ws.Range("L109").Value = 30213.75
$ws.Range("N109").Value = -32293.75
$ws.Range("H123").Value = 42604
$ws.Range("J123").Value = 42604
$ws.Range("L123").Value = 42604
$ws.Range("N123").Value = -47504
$ws.Range("H126").Value = 4311.5835
$ws.Range("I126").Value = 3213.7778
$ws.Range("J126").Value = 7605
$ws.Range("K126").Value = 9641.3334
$ws.Range("L126").Value = 22815
$ws.Range("M126").Value = -7171.3334
$ws.Range("N126").Value = -27755
$ws.Range("H132").Value = 3263.2
$ws.Range("J132").Value = 6597.6
$ws.Range("L132").Value = 19792.8
$ws.Range("N132").Value = -24852.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2880.2559
$ws.Range("I93").Value = 2995.5833
$ws.Range("J93").Value = 2734.5789
$ws.Range("K93").Value = 2995.5833
$ws.Range("L93").Value = 2734.5789
$ws.Range("M93").Value = -1747.5833
$ws.Range("N93").Value = -5230.5789
$ws.Range("H100").Value = 5850832
$ws.Range("I100").Value = 6947388
$ws.Range("J100").Value = 2533.3333
$ws.Range("K100").Value = 6947388
$ws.Range("L100").Value = 2533.3333
$ws.Range("M100").Value = -6946847
$ws.Range("N100").Value = -3615.3333
$ws.Range("H132").Value = 2754.1133
$ws.Range("I132").Value = 1697.5294
$ws.Range("K132").Value = 5092.5882
$ws.Range("M132").Value = -2562.5882

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 618
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 618
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 618
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -848
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180
$ws.Range("H136").Value = 3535.5588
$ws.Range("I136").Value = 602.89655
$ws.Range("J136").Value = 20545
$ws.Range("K136").Value = 1808.68965
$ws.Range("L136").Value = 61635
$ws.Range("M136").Value = 741.3103499999997
$ws.Range("N136").Value = -66735
